# Add a new diary entry (row 4) to the learning log, "26 syys" (Sep 26),
# documenting boilerplate/particle-sim progress towards rendering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new row's cells in the same order the original author would have
# typed them (this also controls the order new entries land in the shared
# string table): date, content, time, quality notes, code notes, hours.
$ws.Range("A4").Value = "26 syys"
$ws.Range("C4").Value = "Kirjan sivut 36-54"
$ws.Range("B4").Value = "12.00-13:30, 13:45-14:15"
$ws.Range("E4").Value = "Pieniä viilauksia koodiin ja projektirakenteeseen esim headerit, ja nimiavaruuksien kertailua sekä konstruktorien otsaketiedosto sekä toteuttava luokka"
$ws.Range("D4").Value = "Huonot yöunet verottivat mutta välillä kun nousi ylös niin pysyi hereillä. Motivaatio korkealla edelleen saada jotain ruudulle myös näkymään, pientä lämmittelyä nämä pari päivää."
$ws.Range("G4").Value = 2

# Match the existing formatting used by row 3: time format + wrap text on
# B (kello), and wrap text on D/E (laatu / huomiot koodista).
$ws.Range("B4").WrapText = $true
$ws.Range("B4").NumberFormat = "h:mm"
$ws.Range("D4").WrapText = $true
$ws.Range("E4").WrapText = $true

# Give the new row the taller height needed to show the longer notes.
$ws.Rows.Item(4).RowHeight = 78

Write-Output "Added diary entry for 26 syys"
